# NIT-9001559939.xlsx update
#  - Refresh "VALOR MORA" total and worker/period counters
#  - Re-sort the Pablo Elias Rodriguez Gomez period rows ascending (1607 -> 2507)
#  - Append a new period row (2508) for Pablo Elias Rodriguez Gomez, taking the
#    "last row" border styling that used to belong to the removed employee
#  - Remove the second worker (German Enrique Brito Atencio) entirely
#  - Column D (Nombre Trabajador) narrows now that the longest name was removed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header summary values
# ---------------------------------------------------------------------------
$ws.Range("E11").Value2 = 3265556   # VALOR MORA total
$ws.Range("C13").Value2 = 1         # Cant. Trabajadores (was 2, German removed)
$ws.Range("F13").Value2 = 110       # Cant. Periodos (was 109, +1 new period)

# ---------------------------------------------------------------------------
# 2. Reverse the 109 existing period rows for Pablo (rows 16-124) so the
#    periods run ascending (1607 .. 2507) instead of descending (2507 .. 1607)
# ---------------------------------------------------------------------------
$firstRow = 16
$lastRow = 124
$nRows = $lastRow - $firstRow + 1
$nCols = 6   # columns B..G

$srcRange = $ws.Range("B$($firstRow):G$($lastRow)")
$data = $srcRange.Value2

$reversed = New-Object 'object[,]' $nRows, $nCols
for ($i = 1; $i -le $nRows; $i++) {
    for ($j = 1; $j -le $nCols; $j++) {
        $reversed[$i - 1, $j - 1] = $data[$nRows - $i + 1, $j]
    }
}
$srcRange.Value2 = $reversed

# ---------------------------------------------------------------------------
# 3. Turn the old "second worker" block (rows 125-126, German Enrique Brito
#    Atencio) into a single new period row (2508) for Pablo, reusing the
#    bottom-border "last row" formatting that row 126 had.
# ---------------------------------------------------------------------------
$ws.Range("B126:J126").Copy() | Out-Null
$ws.Range("B125:J125").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Drop German's second row entirely - shifts every following row up by one
$ws.Rows("126").Delete()

# New trailing period row for Pablo Elias Rodriguez Gomez
$ws.Range("B125").Value2 = "CC"
$ws.Range("C125").Value2 = "13830102"
$ws.Range("D125").Value2 = "PABLO ELIAS RODRIGUEZ GOMEZ"
$ws.Range("E125").Value2 = "2508"
$ws.Range("F125").Value2 = 31249
$ws.Range("G125").Value2 = 781242
$ws.Range("H125").Value2 = $null
$ws.Range("I125").Value2 = $null
$ws.Range("J125").Value2 = $null

# ---------------------------------------------------------------------------
# 4. Column D (Nombre Trabajador) is slightly narrower now that "GERMAN
#    ENRIQUE BRITO ATENCIO" (the longest name) is gone.
# ---------------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 31.333333

Write-Host "Done"
